# PREI / SAI date ranges: add a new "01/01/2025 - 19/06/2025" row above the
# existing "20/06/2025 - 25/08/2025" row, and widen column B so the longer
# second-column values (now used for a date-file-manager export) fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the current data row (row 2) down to row 3, leaving a fresh row 2.
$ws.Rows.Item(2).Insert(-4121)   # xlShiftDown

$ws.Range("A2").Value = "01/01/2025"
$ws.Range("B2").Value = "19/06/2025"

# Row 2 inherited column A's formatting from the header row on insert;
# re-apply the wrapped/text style so A2 matches A3 (and the rest of column A).
$ws.Range("A2").WrapText = $true

# Column B no longer auto-sizes to its contents; give it a fixed, wider width.
$ws.Columns.Item(2).ColumnWidth = 23.83
